# =====================================================================================
# Update countries & provincias Spain
#
# The underlying "Pais" country master list was refreshed upstream: several countries were
# re-ordered in the shared master list (e.g. "Chequia" now sorts before "Irlanda", "Irak"
# before "Hong Kong", etc.) and the statistics snapshot was refreshed from 17:22 to 17:52.
# Each worksheet row keeps its original position, but for the rows whose underlying country
# changed, the country name (column A) and/or the refreshed statistics (columns B:H) need to
# be written with their new values.
# =====================================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------------------------
# 1) Country names (column A) - including the "Datos actualizados..." timestamp in row 1 -
#    that now resolve to a different entry of the re-ordered master list.
# ---------------------------------------------------------------------------------------
$countryNameUpdates = @(
    @{ Row = 1; Name = "Datos actualizados a 4 de Abril de 2020 a las 17:52" }
    @{ Row = 26; Name = "Chequia" }
    @{ Row = 27; Name = "Irlanda" }
    @{ Row = 28; Name = "Chile" }
    @{ Row = 29; Name = "Dinamarca" }
    @{ Row = 30; Name = "Polonia" }
    @{ Row = 31; Name = "Rumania" }
    @{ Row = 63; Name = "Irak" }
    @{ Row = 64; Name = "Hong Kong" }
    @{ Row = 65; Name = "Marruecos" }
    @{ Row = 85; Name = "Republica de Chipre" }
    @{ Row = 86; Name = "Costa Rica" }
    @{ Row = 91; Name = "Burkina Faso" }
    @{ Row = 92; Name = "Jordania" }
    @{ Row = 94; Name = "Cuba" }
    @{ Row = 95; Name = "Oman" }
    @{ Row = 177; Name = "Laos" }
    @{ Row = 178; Name = "Liberia" }
    @{ Row = 179; Name = "Mozambique" }
    @{ Row = 202; Name = "Sierra Leona" }
    @{ Row = 203; Name = "Botsuana" }
    @{ Row = 204; Name = "Gambia" }
    @{ Row = 205; Name = "Islas Virgenes Britanicas" }
    @{ Row = 206; Name = "Anguila" }
    @{ Row = 207; Name = "Burundi" }
    @{ Row = 208; Name = "Bonaire, San Eustaquio y Saba" }
    @{ Row = 209; Name = "Papua Nueva Guinea" }
    @{ Row = 210; Name = "Timor Oriental" }
    @{ Row = 211; Name = "Islas Malvinas" }
)

foreach ($u in $countryNameUpdates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Name
}

# ---------------------------------------------------------------------------------------
# 2) Refreshed statistics (columns B:H = Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose numbers changed
#    between the 17:22 and 17:52 snapshots.
# ---------------------------------------------------------------------------------------
$statUpdates = @(
    @{ Row = 4; Col = 2; Value = 291021 }
    @{ Row = 4; Col = 3; Value = 13860 }
    @{ Row = 4; Col = 4; Value = 14368 }
    @{ Row = 4; Col = 5; Value = 268809 }
    @{ Row = 15; Col = 6; Value = 1360 }
    @{ Row = 16; Col = 5; Value = 10013 }
    @{ Row = 16; Col = 7; Value = 6 }
    @{ Row = 16; Col = 8; Value = 214 }
    @{ Row = 26; Col = 2; Value = 4362 }
    @{ Row = 26; Col = 3; Value = 172 }
    @{ Row = 26; Col = 4; Value = 78 }
    @{ Row = 26; Col = 5; Value = 4225 }
    @{ Row = 26; Col = 6; Value = 87 }
    @{ Row = 26; Col = 7; Value = 6 }
    @{ Row = 26; Col = 8; Value = 59 }
    @{ Row = 27; Col = 2; Value = 4273 }
    @{ Row = 27; Col = 3; Value = 0 }
    @{ Row = 27; Col = 4; Value = 25 }
    @{ Row = 27; Col = 5; Value = 4128 }
    @{ Row = 27; Col = 6; Value = 148 }
    @{ Row = 27; Col = 7; Value = 0 }
    @{ Row = 27; Col = 8; Value = 120 }
    @{ Row = 28; Col = 2; Value = 4161 }
    @{ Row = 28; Col = 3; Value = 424 }
    @{ Row = 28; Col = 4; Value = 427 }
    @{ Row = 28; Col = 5; Value = 3707 }
    @{ Row = 28; Col = 6; Value = 31 }
    @{ Row = 28; Col = 7; Value = 5 }
    @{ Row = 28; Col = 8; Value = 27 }
    @{ Row = 29; Col = 2; Value = 4077 }
    @{ Row = 29; Col = 3; Value = 320 }
    @{ Row = 29; Col = 4; Value = 1283 }
    @{ Row = 29; Col = 5; Value = 2633 }
    @{ Row = 29; Col = 6; Value = 142 }
    @{ Row = 29; Col = 7; Value = 22 }
    @{ Row = 29; Col = 8; Value = 161 }
    @{ Row = 30; Col = 2; Value = 3627 }
    @{ Row = 30; Col = 3; Value = 244 }
    @{ Row = 30; Col = 4; Value = 116 }
    @{ Row = 30; Col = 5; Value = 3432 }
    @{ Row = 30; Col = 6; Value = 50 }
    @{ Row = 30; Col = 8; Value = 79 }
    @{ Row = 31; Col = 2; Value = 3613 }
    @{ Row = 31; Col = 3; Value = 430 }
    @{ Row = 31; Col = 4; Value = 329 }
    @{ Row = 31; Col = 5; Value = 3143 }
    @{ Row = 31; Col = 6; Value = 119 }
    @{ Row = 31; Col = 7; Value = 8 }
    @{ Row = 31; Col = 8; Value = 141 }
    @{ Row = 37; Col = 2; Value = 2724 }
    @{ Row = 37; Col = 3; Value = 38 }
    @{ Row = 37; Col = 5; Value = 2553 }
    @{ Row = 59; Col = 2; Value = 1039 }
    @{ Row = 59; Col = 3; Value = 78 }
    @{ Row = 59; Col = 5; Value = 967 }
    @{ Row = 63; Col = 2; Value = 878 }
    @{ Row = 63; Col = 3; Value = 58 }
    @{ Row = 63; Col = 4; Value = 259 }
    @{ Row = 63; Col = 5; Value = 563 }
    @{ Row = 63; Col = 6; Value = 0 }
    @{ Row = 63; Col = 7; Value = 2 }
    @{ Row = 63; Col = 8; Value = 56 }
    @{ Row = 64; Col = 2; Value = 862 }
    @{ Row = 64; Col = 3; Value = 17 }
    @{ Row = 64; Col = 4; Value = 173 }
    @{ Row = 64; Col = 5; Value = 685 }
    @{ Row = 64; Col = 6; Value = 8 }
    @{ Row = 64; Col = 7; Value = 0 }
    @{ Row = 64; Col = 8; Value = 4 }
    @{ Row = 65; Col = 2; Value = 858 }
    @{ Row = 65; Col = 3; Value = 67 }
    @{ Row = 65; Col = 4; Value = 62 }
    @{ Row = 65; Col = 5; Value = 746 }
    @{ Row = 65; Col = 6; Value = 1 }
    @{ Row = 65; Col = 7; Value = 2 }
    @{ Row = 65; Col = 8; Value = 50 }
    @{ Row = 75; Col = 4; Value = 54 }
    @{ Row = 75; Col = 5; Value = 449 }
    @{ Row = 85; Col = 2; Value = 426 }
    @{ Row = 85; Col = 3; Value = 30 }
    @{ Row = 85; Col = 4; Value = 28 }
    @{ Row = 85; Col = 5; Value = 387 }
    @{ Row = 85; Col = 6; Value = 11 }
    @{ Row = 85; Col = 8; Value = 11 }
    @{ Row = 86; Col = 2; Value = 416 }
    @{ Row = 86; Col = 4; Value = 11 }
    @{ Row = 86; Col = 5; Value = 403 }
    @{ Row = 86; Col = 6; Value = 13 }
    @{ Row = 86; Col = 8; Value = 2 }
    @{ Row = 91; Col = 2; Value = 318 }
    @{ Row = 91; Col = 3; Value = 16 }
    @{ Row = 91; Col = 4; Value = 66 }
    @{ Row = 91; Col = 5; Value = 236 }
    @{ Row = 91; Col = 6; Value = 0 }
    @{ Row = 91; Col = 8; Value = 16 }
    @{ Row = 92; Col = 2; Value = 310 }
    @{ Row = 92; Col = 4; Value = 58 }
    @{ Row = 92; Col = 5; Value = 247 }
    @{ Row = 92; Col = 6; Value = 5 }
    @{ Row = 92; Col = 8; Value = 5 }
    @{ Row = 94; Col = 2; Value = 288 }
    @{ Row = 94; Col = 3; Value = 19 }
    @{ Row = 94; Col = 4; Value = 15 }
    @{ Row = 94; Col = 5; Value = 267 }
    @{ Row = 94; Col = 6; Value = 8 }
    @{ Row = 94; Col = 8; Value = 6 }
    @{ Row = 95; Col = 2; Value = 277 }
    @{ Row = 95; Col = 3; Value = 25 }
    @{ Row = 95; Col = 4; Value = 61 }
    @{ Row = 95; Col = 5; Value = 215 }
    @{ Row = 95; Col = 6; Value = 3 }
    @{ Row = 95; Col = 8; Value = 1 }
    @{ Row = 127; Col = 2; Value = 77 }
    @{ Row = 127; Col = 3; Value = 2 }
    @{ Row = 127; Col = 5; Value = 77 }
    @{ Row = 178; Col = 3; Value = 3 }
    @{ Row = 178; Col = 5; Value = 9 }
    @{ Row = 178; Col = 7; Value = 1 }
    @{ Row = 178; Col = 8; Value = 1 }
    @{ Row = 179; Col = 3; Value = 0 }
    @{ Row = 179; Col = 4; Value = 1 }
    @{ Row = 179; Col = 7; Value = 0 }
    @{ Row = 179; Col = 8; Value = 0 }
    @{ Row = 202; Col = 3; Value = 2 }
    @{ Row = 202; Col = 5; Value = 4 }
    @{ Row = 202; Col = 8; Value = 0 }
    @{ Row = 203; Col = 4; Value = 0 }
    @{ Row = 203; Col = 5; Value = 3 }
    @{ Row = 204; Col = 2; Value = 4 }
    @{ Row = 204; Col = 4; Value = 2 }
    @{ Row = 204; Col = 5; Value = 1 }
    @{ Row = 204; Col = 8; Value = 1 }
    @{ Row = 207; Col = 2; Value = 3 }
    @{ Row = 207; Col = 5; Value = 3 }
)

foreach ($u in $statUpdates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
